# Atualização de bases das ligas, do dia: 18-04-2024 às 00:36
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 179 and 180 swap their match-data (B, F..AC); columns A, C, D, E
#    (id / Div / Div Original Name / Date) stay put since they key the row.
# ---------------------------------------------------------------------------

# --- new content for row 179 (previously row 180's data) ---
$ws.Range("B179").Value = 8026714
$ws.Range("F179").Value = "BG Pathum United"
$ws.Range("G179").Value = "Buriram United"
$ws.Range("H179").Value = 1
$ws.Range("I179").Value = 1
$ws.Range("J179").Value = "D"
$ws.Range("K179").Value = 3
$ws.Range("L179").Value = 3.6
$ws.Range("M179").Value = 2
$ws.Range("N179").Value = 3.1
$ws.Range("O179").Value = 3.75
$ws.Range("P179").Value = 1.95
$ws.Range("Q179").Value = 0.5
$ws.Range("R179").Value = 1.825
$ws.Range("S179").Value = 1.975
$ws.Range("T179").Value = 2.75
$ws.Range("U179").Value = 1.85
$ws.Range("V179").Value = 1.95
$ws.Range("W179").Value = -1
$ws.Range("X179").Value = 2.75
$ws.Range("Y179").Value = -1
$ws.Range("Z179").Value = 0.825
$ws.Range("AA179").Value = -1
$ws.Range("AB179").Value = -1
$ws.Range("AC179").Value = 0.95

# --- new content for row 180 (previously row 179's data) ---
$ws.Range("B180").Value = 6992695
$ws.Range("F180").Value = "Muang Thong United"
$ws.Range("G180").Value = "Uthai Thani FC"
$ws.Range("H180").Value = 5
$ws.Range("I180").Value = 2
$ws.Range("J180").Value = "H"
$ws.Range("K180").Value = 2.1
$ws.Range("L180").Value = 3.75
$ws.Range("M180").Value = 2.7
$ws.Range("N180").Value = 1.95
$ws.Range("O180").Value = 3.8
$ws.Range("P180").Value = 2.9
$ws.Range("Q180").Value = -0.25
$ws.Range("R180").Value = 1.8
$ws.Range("S180").Value = 2
$ws.Range("T180").Value = 3
$ws.Range("U180").Value = 1.825
$ws.Range("V180").Value = 1.975
$ws.Range("W180").Value = 0.95
$ws.Range("X180").Value = -1
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = 0.8
$ws.Range("AA180").Value = -1
$ws.Range("AB180").Value = 0.825
$ws.Range("AC180").Value = -1

# ---------------------------------------------------------------------------
# 2) Three new fixtures appended at the bottom of the sheet (rows 191-193).
#    H/I/J (score/result) are left blank - matches not played yet.
#    A (id) and E (Date) reuse the same cell styles as the row above them.
# ---------------------------------------------------------------------------

# --- row 191 : id 189 ---
$ws.Range("A191").Value = 189
$ws.Range("B191").Value = 8075058
$ws.Range("C191").Value = "Thailand Premier League"
$ws.Range("D191").Value = "Thailand Premier League"
$ws.Range("E191").Value = 45400.375
$ws.Range("F191").Value = "Bangkok United"
$ws.Range("G191").Value = "Lamphun Warrior FC"
$ws.Range("K191").Value = 1.363
$ws.Range("L191").Value = 4.333
$ws.Range("M191").Value = 6.5
$ws.Range("N191").Value = 1.4
$ws.Range("O191").Value = 4.2
$ws.Range("P191").Value = 6
$ws.Range("Q191").Value = -1.25
$ws.Range("R191").Value = 1.875
$ws.Range("S191").Value = 1.925
$ws.Range("T191").Value = 2.75
$ws.Range("U191").Value = 1.8
$ws.Range("V191").Value = 2
$ws.Range("W191").Value = 0
$ws.Range("X191").Value = 0
$ws.Range("Y191").Value = 0
$ws.Range("Z191").Value = 0
$ws.Range("AA191").Value = 0

# --- row 192 : id 190 ---
$ws.Range("A192").Value = 190
$ws.Range("B192").Value = 6992709
$ws.Range("C192").Value = "Thailand Premier League"
$ws.Range("D192").Value = "Thailand Premier League"
$ws.Range("E192").Value = 45401.375
$ws.Range("F192").Value = "Port FC"
$ws.Range("G192").Value = "Chiangrai Utd"
$ws.Range("K192").Value = 1.363
$ws.Range("L192").Value = 4.5
$ws.Range("M192").Value = 7.5
$ws.Range("N192").Value = 1.285
$ws.Range("O192").Value = 5
$ws.Range("P192").Value = 8.5
$ws.Range("Q192").Value = -1.5
$ws.Range("R192").Value = 1.85
$ws.Range("S192").Value = 1.95
$ws.Range("T192").Value = 2.75
$ws.Range("U192").Value = 1.8
$ws.Range("V192").Value = 2
$ws.Range("W192").Value = 0
$ws.Range("X192").Value = 0
$ws.Range("Y192").Value = 0
$ws.Range("Z192").Value = 0
$ws.Range("AA192").Value = 0

# --- row 193 : id 191 ---
$ws.Range("A193").Value = 191
$ws.Range("B193").Value = 6995900
$ws.Range("C193").Value = "Thailand Premier League"
$ws.Range("D193").Value = "Thailand Premier League"
$ws.Range("E193").Value = 45402.35416666666
$ws.Range("F193").Value = "Police Tero FC"
$ws.Range("G193").Value = "Uthai Thani FC"
$ws.Range("K193").Value = 3.3
$ws.Range("L193").Value = 3.6
$ws.Range("M193").Value = 1.95
$ws.Range("N193").Value = 3.3
$ws.Range("O193").Value = 3.6
$ws.Range("P193").Value = 1.95
$ws.Range("Q193").Value = 0.5
$ws.Range("R193").Value = 1.8
$ws.Range("S193").Value = 2
$ws.Range("T193").Value = 3
$ws.Range("U193").Value = 1.975
$ws.Range("V193").Value = 1.825
$ws.Range("W193").Value = 0
$ws.Range("X193").Value = 0
$ws.Range("Y193").Value = 0
$ws.Range("Z193").Value = 0
$ws.Range("AA193").Value = 0

# ---------------------------------------------------------------------------
# 3) Re-apply the id (bold+border) and Date (date-number-format) cell styles
#    to the new rows by copying them (format only) from the row above.
# ---------------------------------------------------------------------------
$ws.Range("A190").Copy()
$ws.Range("A191").PasteSpecial(-4122)
$ws.Range("A191").Copy()
$ws.Range("A192").PasteSpecial(-4122)
$ws.Range("A192").Copy()
$ws.Range("A193").PasteSpecial(-4122)

$ws.Range("E190").Copy()
$ws.Range("E191").PasteSpecial(-4122)
$ws.Range("E191").Copy()
$ws.Range("E192").PasteSpecial(-4122)
$ws.Range("E192").Copy()
$ws.Range("E193").PasteSpecial(-4122)

$excel.CutCopyMode = 0
